$d = $word.ActiveDocument

$replacements = @(
    @("85×35=2975", "28×67=1876"),
    @("94×48=4512", "81×54=4374"),
    @("49×59=2891", "74×46=3404"),
    @("97×37=3589", "88×57=5016"),
    @("47×84=3948", "58×88=5104"),
    @("18×95=1710", "36×59=2124"),
    @("61×19=1159", "53×21=1113"),
    @("15×91=1365", "12×37=444"),
    @("67×43=2881", "81×63=5103"),
    @("84×96=8064", "63×53=3339"),
    @("94×85=7990", "11×73=803"),
    @("13×55=715", "53×36=1908"),
    @("54×89=4806", "70×60=4200"),
    @("58×91=5278", "30×47=1410"),
    @("29×62=1798", "81×62=5022"),
    @("44×64=2816", "91×65=5915"),
    @("79×95=7505", "77×26=2002"),
    @("99×79=7821", "53×26=1378"),
    @("28×90=2520", "68×93=6324"),
    @("22×36=792",  "57×42=2394"),
    @("79×12=948",  "66×79=5214"),
    @("93×64=5952", "79×31=2449"),
    @("60×47=2820", "72×55=3960"),
    @("19×39=741",  "16×12=192"),
    @("66×45=2970", "68×28=1904")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
